$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.5022738159769057
$ws.Range("J2").Value = 0.5022738159769057
$ws.Range("O2").Value = 0.06522509891308133
$ws.Range("P2").Value = 0.06522509891308133
$ws.Range("S2").Value = 0.03276085932854448
$ws.Range("T2").Value = 0.03276085932854448
$ws.Range("I3").Value = 0.5022738159769057
$ws.Range("J3").Value = 0.5022738159769057
$ws.Range("M3").Value = 0.6481333333333333
$ws.Range("N3").Value = 1.9444
$ws.Range("O3").Value = 0.2064033004146749
$ws.Range("P3").Value = 0.2064033004146749
$ws.Range("Q3").Value = 0.4131862962666667
$ws.Range("R3").Value = 3.7186766664
$ws.Range("S3").Value = 0.1036709733295064
$ws.Range("T3").Value = 0.1036709733295064
$ws.Range("I4").Value = 0.5022738159769057
$ws.Range("J4").Value = 0.5022738159769057
$ws.Range("M4").Value = 1.888205
$ws.Range("N4").Value = 5.664615
$ws.Range("O4").Value = 0.6013141491351952
$ws.Range("P4").Value = 0.6013141491351952
$ws.Range("Q4").Value = 1.20373446391
$ws.Range("R4").Value = 10.83361017519
$ws.Range("S4").Value = 0.3020243522870407
$ws.Range("T4").Value = 0.3020243522870407
$ws.Range("I5").Value = 0.5022738159769057
$ws.Range("J5").Value = 0.5022738159769057
$ws.Range("M5").Value = 0.398977
$ws.Range("N5").Value = 1.196931
$ws.Range("O5").Value = 0.1270574515370486
$ws.Range("P5").Value = 0.1270574515370486
$ws.Range("Q5").Value = 0.254348635454
$ws.Range("R5").Value = 2.289137719086
$ws.Range("S5").Value = 0.06381763103181413
$ws.Range("T5").Value = 0.06381763103181413
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.63173
$ws.Range("H6").Value = 1.89519
$ws.Range("I6").Value = 0.4977261840230943
$ws.Range("J6").Value = 0.4977261840230943
$ws.Range("O6").Value = 0.06522509891308133
$ws.Range("P6").Value = 0.06522509891308133
$ws.Range("Q6").Value = 0.1293879905266667
$ws.Range("R6").Value = 1.16449191474
$ws.Range("S6").Value = 0.03246423958453684
$ws.Range("T6").Value = 0.03246423958453685
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.63173
$ws.Range("H7").Value = 1.89519
$ws.Range("I7").Value = 0.4977261840230943
$ws.Range("J7").Value = 0.4977261840230943
$ws.Range("M7").Value = 0.6481333333333333
$ws.Range("N7").Value = 1.9444
$ws.Range("O7").Value = 0.2064033004146749
$ws.Range("P7").Value = 0.2064033004146749
$ws.Range("Q7").Value = 0.4094452706666667
$ws.Range("R7").Value = 3.685007436
$ws.Range("S7").Value = 0.1027323270851685
$ws.Range("T7").Value = 0.1027323270851685
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.63173
$ws.Range("H8").Value = 1.89519
$ws.Range("I8").Value = 0.4977261840230943
$ws.Range("J8").Value = 0.4977261840230943
$ws.Range("M8").Value = 1.888205
$ws.Range("N8").Value = 5.664615
$ws.Range("O8").Value = 0.6013141491351952
$ws.Range("P8").Value = 0.6013141491351952
$ws.Range("Q8").Value = 1.19283574465
$ws.Range("R8").Value = 10.73552170185
$ws.Range("S8").Value = 0.2992897968481545
$ws.Range("T8").Value = 0.2992897968481545
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.63173
$ws.Range("H9").Value = 1.89519
$ws.Range("I9").Value = 0.4977261840230943
$ws.Range("J9").Value = 0.4977261840230943
$ws.Range("M9").Value = 0.398977
$ws.Range("N9").Value = 1.196931
$ws.Range("O9").Value = 0.1270574515370486
$ws.Range("P9").Value = 0.1270574515370486
$ws.Range("Q9").Value = 0.25204574021
$ws.Range("R9").Value = 2.26841166189
$ws.Range("S9").Value = 0.06323982050523441
$ws.Range("T9").Value = 0.06323982050523443